$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "59.064.45"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "2.575.52"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "508.35"
$ws.Range("E5").Value = "  -2.88%  "
Set-TextValue $ws.Range("D6") "143.73"
$ws.Range("E6").Value = "  -7.11%  "
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -5.76%  "
$ws.Range("D9").Value = "2.587.81"
$ws.Range("E9").Value = "  -0.82%  "
Set-TextValue $ws.Range("D10") "6.23"
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "3.028.77"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "59.057.38"
$ws.Range("E15").Value = "  -3.11%  "
Set-TextValue $ws.Range("D16") "20.65"
$ws.Range("E16").Value = "  -4.68%  "
$ws.Range("E17").Value = "  -4.47%  "
$ws.Range("D18").Value = "2.577.78"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  -5.02%  "
Set-TextValue $ws.Range("D20") "333.54"
$ws.Range("E20").Value = "  -6.07%  "
$ws.Range("E21").Value = "  -4.63%  "
Set-TextValue $ws.Range("D22") "1.00"
$ws.Range("E22").Value = "  +0.11%  "
Set-TextValue $ws.Range("D23") "5.97"
$ws.Range("E23").Value = "  -3.94%  "
Set-TextValue $ws.Range("D24") "59.80"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").Value = "  -4.27%  "
Set-TextValue $ws.Range("D26") "0.997"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -5.47%  "
$ws.Range("E28").Value = "  -7.74%  "
Set-TextValue $ws.Range("D29") "6.89"
$ws.Range("E29").Value = "  -7.03%  "
Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D31") "18.62"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D32") "149.63"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D33") "5.84"
$ws.Range("E33").Value = "  -6.58%  "
$ws.Range("E34").Value = "  -3.56%  "
Set-TextValue $ws.Range("D35") "3.95"
$ws.Range("E35").Value = "  -5.70%  "
Set-TextValue $ws.Range("D36") "0.895"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("E37").Value = "  -7.99%  "
$ws.Range("E38").Value = "  -1.52%  "
Set-TextValue $ws.Range("D39") "0.826"
$ws.Range("E39").Value = "  -7.33%  "
Set-TextValue $ws.Range("D40") "288.63"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  -7.50%  "
Set-TextValue $ws.Range("D42") "3.51"
$ws.Range("E42").Value = "  -7.90%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.612"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "0.998"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -3.28%  "
Set-TextValue $ws.Range("D46") "0.0532"
$ws.Range("E46").Value = "  -4.95%  "
Set-TextValue $ws.Range("D47") "18.78"
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  -4.54%  "
Set-TextValue $ws.Range("D50") "4.54"
$ws.Range("E50").Value = "  -7.55%  "
$ws.Range("D51").Value = "1.919.60"
$ws.Range("E51").Value = "  -2.16%  "
